$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text for the d9c72ebc file changes from "Ready for handoff" to
#    "Handback transform failed" everywhere it is shown (Overview!E3/F3,
#    zh-cn!C3, de-de!C3 all share the same underlying text).
# ---------------------------------------------------------------------------
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value     = "Handback transform failed"
$wsDeDe.Range("C3").Value     = "Handback transform failed"

# ---------------------------------------------------------------------------
# 2. New "Error Detail" messages (column P) explaining the handback
#    transform failure for the d9c72ebc file, one per locale.
# ---------------------------------------------------------------------------
$wsZhCn.Range("P3").Value = "Handback file name: 3u12fzi5.xni is different with handoff file name: d9c72ebc-bf30-4cc5-aa35-9fcc7ea7d485.57ec2417f2a31064dff290bb6fc197e75560e4db.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: 3u12fzi5.xni is different with handoff file name: d9c72ebc-bf30-4cc5-aa35-9fcc7ea7d485.57ec2417f2a31064dff290bb6fc197e75560e4db.de-de."

# ---------------------------------------------------------------------------
# 3. Widen column P ("Error Detail") on both locale sheets so the new
#    message text is readable.
# ---------------------------------------------------------------------------
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
